$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibitions)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 90
$ws1.Range("F4").Value = 400
$ws1.Range("F5").Value = 180
$ws1.Range("F6").Value = 126
$ws1.Range("F7").Value = 1058
$ws1.Range("G7").Value = 61.2
$ws1.Range("F8").Value = 351
$ws1.Range("F9").Value = 179
$ws1.Range("F10").Value = 46
$ws1.Range("F13").Value = 356
$ws1.Range("F14").Value = 767
$ws1.Range("F15").Value = 139
$ws1.Range("F16").Value = 704
$ws1.Range("F17").Value = 266
$ws1.Range("F18").Value = 68
$ws1.Range("F19").Value = 981
$ws1.Range("F20").Value = 440
$ws1.Range("F21").Value = 252
$ws1.Range("F22").Value = 79
$ws1.Range("F23").Value = 365
$ws1.Range("F25").Value = 36
$ws1.Range("F26").Value = 460

# Sheet 2: 演出 (Performances)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F4").Value = 358
$ws2.Range("F6").Value = 41
$ws2.Range("F13").Value = 13

# Sheet 4: 全部类型 (All Types)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F4").Value = 90
$ws4.Range("F6").Value = 400
$ws4.Range("F7").Value = 180
$ws4.Range("F8").Value = 126
$ws4.Range("F9").Value = 1058
$ws4.Range("G9").Value = 61.2
$ws4.Range("F10").Value = 351
$ws4.Range("F11").Value = 179
$ws4.Range("F13").Value = 46
$ws4.Range("F14").Value = 358
$ws4.Range("F18").Value = 41
$ws4.Range("F20").Value = 356
$ws4.Range("F21").Value = 767
$ws4.Range("F22").Value = 139
$ws4.Range("F23").Value = 704
$ws4.Range("F24").Value = 266
$ws4.Range("F25").Value = 68
$ws4.Range("F26").Value = 981
$ws4.Range("F27").Value = 440
$ws4.Range("F30").Value = 252
$ws4.Range("F31").Value = 79
$ws4.Range("F32").Value = 365
$ws4.Range("F36").Value = 36
$ws4.Range("F38").Value = 460
$ws4.Range("F39").Value = 13
